$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (Behavior) text updates, issued in the exact order the new
# shared-string entries should be appended to the shared strings table ---
$ws.Range("C2").Value = "TANG"
$ws.Range("C12").Value = "BLOCKY"
$ws.Range("C13").Value = "BLOCKN"
$ws.Range("C16").Value = "EMESIS"
$ws.Range("C36").Value = "NO DATA"
$ws.Range("C43").Value = "PACK"
$ws.Range("C10").Value = "JAW P"
$ws.Range("C11").Value = "CHIN P"
$ws.Range("C8").Value = "FING P"
$ws.Range("C7").Value = "PHYS G"
$ws.Range("C21").Value = "RP SOL"
$ws.Range("C22").Value = "RP SACC"
$ws.Range("C23").Value = "RP LIQ"
$ws.Range("C24").Value = "RP LACC"
$ws.Range("C26").Value = "LIQ ACC"
$ws.Range("C42").Value = "LIQ PRES"
$ws.Range("C28").Value = "SOL ACC"
$ws.Range("C29").Value = "SOL PRES"
$ws.Range("C30").Value = "SOL ABS"
$ws.Range("C31").Value = "SOL PRO"
$ws.Range("C37").Value = "IND LIQ"
$ws.Range("C38").Value = "IND SOL"
$ws.Range("C39").Value = "LIQ MC"
$ws.Range("C40").Value = "SOL MC"
$ws.Range("C41").Value = "LIQ ABS"
$ws.Range("C27").Value = "LIQ PRO"
$ws.Range("C35").Value = "DUMMY 1"
$ws.Range("C6").Value = "DUM DUR"
$ws.Range("C46").Value = "CORRECT"

# --- Column B (Type) updates -- these reuse existing shared strings ---
$ws.Range("B6").Value = "Duration"
$ws.Range("B17").Value = "Duration"
$ws.Range("B18").Value = "Duration"
$ws.Range("B27").Value = "Once"
$ws.Range("B35").Value = "Once"
$ws.Range("B42").Value = "Duration"
$ws.Range("B46").Value = "Duration"

# --- Update the saved selection / scroll position to match the author's last view ---
$ws.Range("D28").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
